# Generate Report for Handoff
# The "b.md" file has now been handed off: update its status from
# "Handed back: in sync with en-US" to "Ready for handoff" on all sheets,
# and record the newly generated handoff file names / timestamps on the
# zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 is the "b.md" row ---
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: row 3 is the "b.md" row ---
$ws2.Range("B3").Value = "Ready for handoff"
$ws2.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$ws2.Range("D3").Value = "2016-02-18 07:47:45"

foreach ($hl in $ws2.Hyperlinks) {
    if ($hl.Range.Address() -eq '$C$3') {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
    }
}

# --- de-de sheet: row 3 is the "b.md" row ---
$ws3.Range("B3").Value = "Ready for handoff"
$ws3.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$ws3.Range("D3").Value = "2016-02-18 07:48:00"

foreach ($hl in $ws3.Hyperlinks) {
    if ($hl.Range.Address() -eq '$C$3') {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
    }
}
